# Trade #15 closed at 2026-02-17 13:17:31 - unknown UNKNOWN +0.000%
#
# Updates the Summary and Strategy Status aggregate figures to reflect the
# newly-closed trade, and appends the new trade row (#15, zero-indexed
# Trade # = 15) to both the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# --- Summary sheet: refresh aggregate stats ---------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.25   # Current Capital
$summary.Range("B4").Value = -0.75     # Total P&L $
$summary.Range("B5").Value = -1        # Total P&L %
$summary.Range("B6").Value = 15        # Total Trades
$summary.Range("B8").Value = 10        # Losing Trades
$summary.Range("B9").Value = 26.67     # Win Rate %

# --- Strategy Status sheet: refresh MarketMaking row (row 4) ----------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.25      # Capital
$status.Range("D4").Value = 15         # Trades
$status.Range("E4").Value = -0.75      # P&L $
$status.Range("F4").Value = -0.75      # P&L %
$status.Range("G4").Value = 26.67      # Win Rate %

# --- Append the new closed trade (row 16) to both trade logs ----------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(16, 1).Value = 15

    # Force the date to stay plain text instead of being auto-parsed into a
    # date serial number by the COM layer's smart-typing.
    $ws.Cells.Item(16, 2).NumberFormat = "@"
    $ws.Cells.Item(16, 2).Value = "2026-02-17"
    $ws.Cells.Item(16, 2).Style = "Normal"

    $ws.Cells.Item(16, 3).Value = "13:17:25"
    $ws.Cells.Item(16, 4).Value = "MarketMaking"
    $ws.Cells.Item(16, 5).Value = "UP"
    $ws.Cells.Item(16, 6).Value = 0.35
    $ws.Cells.Item(16, 7).Value = 0.23
    $ws.Cells.Item(16, 8).Value = "CLOSED"
    $ws.Cells.Item(16, 9).Value = -34.2857
    $ws.Cells.Item(16, 10).Value = -0.12
    $ws.Cells.Item(16, 11).Value = 99.25
    $ws.Cells.Item(16, 12).Value = 0
    $ws.Cells.Item(16, 13).Value = 0
    $ws.Cells.Item(16, 14).Value = 0.6
    $ws.Cells.Item(16, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(16, 16).Value = "early_exit"
    $ws.Cells.Item(16, 17).Value = 0.14
}
